$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The previous "Estados de Cuenta" (EC) rows are removed and replaced
# with a freshly re-sorted/updated set (same 3 workers x 6 periods,
# but re-grouped by worker and with new "Salario Basico" / "Valor Mora"
# numbers for FABIO ANDRES ARROYO BELTRAN and for period 1903).
# -----------------------------------------------------------------

$rows = @(
    @{ R = 16; Doc = "73009947"; Nombre = "FABIO ANDRES ARROYO BELTRAN";  Periodo = "1903"; Mora = 40000;  Salario = 1200000 },
    @{ R = 17; Doc = "73009947"; Nombre = "FABIO ANDRES ARROYO BELTRAN";  Periodo = "1902"; Mora = 48000;  Salario = 1200000 },
    @{ R = 18; Doc = "73009947"; Nombre = "FABIO ANDRES ARROYO BELTRAN";  Periodo = "1901"; Mora = 48000;  Salario = 1200000 },
    @{ R = 19; Doc = "73009947"; Nombre = "FABIO ANDRES ARROYO BELTRAN";  Periodo = "1812"; Mora = 48000;  Salario = 1200000 },
    @{ R = 20; Doc = "73009947"; Nombre = "FABIO ANDRES ARROYO BELTRAN";  Periodo = "1811"; Mora = 48000;  Salario = 1200000 },
    @{ R = 21; Doc = "73009947"; Nombre = "FABIO ANDRES ARROYO BELTRAN";  Periodo = "1810"; Mora = 48000;  Salario = 1200000 },
    @{ R = 22; Doc = "45478050"; Nombre = "MARYSEL CAÑAS PALACIO";        Periodo = "1903"; Mora = 26667;  Salario = 800000  },
    @{ R = 23; Doc = "45478050"; Nombre = "MARYSEL CAÑAS PALACIO";        Periodo = "1902"; Mora = 32000;  Salario = 800000  },
    @{ R = 24; Doc = "45478050"; Nombre = "MARYSEL CAÑAS PALACIO";        Periodo = "1901"; Mora = 32000;  Salario = 800000  },
    @{ R = 25; Doc = "45478050"; Nombre = "MARYSEL CAÑAS PALACIO";        Periodo = "1812"; Mora = 32000;  Salario = 800000  },
    @{ R = 26; Doc = "45478050"; Nombre = "MARYSEL CAÑAS PALACIO";        Periodo = "1811"; Mora = 32000;  Salario = 800000  },
    @{ R = 27; Doc = "45478050"; Nombre = "MARYSEL CAÑAS PALACIO";        Periodo = "1810"; Mora = 32000;  Salario = 800000  },
    @{ R = 28; Doc = "73572972"; Nombre = "OSCAR LUIS MENDIETA ESTARITA"; Periodo = "1903"; Mora = 50000;  Salario = 1500000 },
    @{ R = 29; Doc = "73572972"; Nombre = "OSCAR LUIS MENDIETA ESTARITA"; Periodo = "1902"; Mora = 60000;  Salario = 1500000 },
    @{ R = 30; Doc = "73572972"; Nombre = "OSCAR LUIS MENDIETA ESTARITA"; Periodo = "1901"; Mora = 60000;  Salario = 1500000 },
    @{ R = 31; Doc = "73572972"; Nombre = "OSCAR LUIS MENDIETA ESTARITA"; Periodo = "1812"; Mora = 60000;  Salario = 1500000 },
    @{ R = 32; Doc = "73572972"; Nombre = "OSCAR LUIS MENDIETA ESTARITA"; Periodo = "1811"; Mora = 60000;  Salario = 1500000 },
    @{ R = 33; Doc = "73572972"; Nombre = "OSCAR LUIS MENDIETA ESTARITA"; Periodo = "1810"; Mora = 60000;  Salario = 1500000 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("C$r").Value = $row.Doc
    $ws.Range("D$r").Value = $row.Nombre
    $ws.Range("E$r").Value = $row.Periodo
    $ws.Range("F$r").Value = $row.Mora
    $ws.Range("G$r").Value = $row.Salario
}

# -----------------------------------------------------------------
# The longer/shorter worker names & document numbers change the
# "best fit" width Excel recalculates for several columns. Reproduce
# the resulting stored column widths as closely as this host allows.
# -----------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
